# Insert a new weekly price record for "Vega Monumental Concepción - Naranja"
# at row 532, pushing every existing record (old rows 532-579) down by one
# row (new rows 533-580).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(532).Insert()

$ws.Cells.Item(532, 1).Value  = 11
$ws.Cells.Item(532, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(532, 3).Value  = "Bíobío"
$ws.Cells.Item(532, 4).Value  = 45265
$ws.Cells.Item(532, 5).Value  = 8
$ws.Cells.Item(532, 6).Value  = "Fruta"
$ws.Cells.Item(532, 7).Value  = 100102
$ws.Cells.Item(532, 8).Value  = "Cítricos"
$ws.Cells.Item(532, 9).Value  = 100102005
$ws.Cells.Item(532, 10).Value = "Naranja"
$ws.Cells.Item(532, 11).Value = "Valencia"
$ws.Cells.Item(532, 12).Value = "Primera"
$ws.Cells.Item(532, 13).Value = 150
$ws.Cells.Item(532, 14).Value = 12000
$ws.Cells.Item(532, 15).Value = 12000
$ws.Cells.Item(532, 16).Value = 12000
$ws.Cells.Item(532, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(532, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(532, 19).Value = 800
$ws.Cells.Item(532, 20).Value = 15
